# Auto-applies the cryptos.xlsx price/volume refresh described in the commit diff.
# D-column "Price" cells are plain text in the original file (t="inlineStr"), even
# when they look like plain numbers (e.g. "1.00", "9.78"). Assigning such a string
# directly via .Value lets Excel auto-coerce it to a real number (losing trailing
# zeros / introducing float noise), so numeric-looking text is written through a
# temporary "@" (text) number format and then the format stamp is cleared again so
# the cell keeps its original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "66.636.14"
$ws.Range("E2").Value = "  -3.53%  "

# Row 3
$ws.Range("D3").Value = "3.560.10"
$ws.Range("E3").Value = "  -4.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
Set-TextCell $ws.Range("D5") "574.37"
$ws.Range("E5").Value = "  -6.62%  "

# Row 6
Set-TextCell $ws.Range("D6") "186.07"
$ws.Range("E6").Value = "  -3.72%  "

# Row 7
$ws.Range("D7").Value = "3.556.29"
$ws.Range("E7").Value = "  -4.35%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.613"
$ws.Range("E8").Value = "  -4.10%  "

# Row 9
$ws.Range("E9").Value = "  +0.32%  "

# Row 10
Set-TextCell $ws.Range("D10") "0.673"
$ws.Range("E10").Value = "  -7.40%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.148"
$ws.Range("E11").Value = "  -8.35%  "

# Row 12
Set-TextCell $ws.Range("D12") "55.08"
$ws.Range("E12").Value = "  -8.27%  "

# Row 13
Set-TextCell $ws.Range("D13") "0.0000261"
$ws.Range("E13").Value = "  -10.30%  "

# Row 14
Set-TextCell $ws.Range("D14") "9.78"
$ws.Range("E14").Value = "  -6.90%  "

# Row 15
$ws.Range("D15").Value = "4.131.43"
$ws.Range("E15").Value = "  -4.23%  "

# Row 16
$ws.Range("D16").Value = "3.558.71"
$ws.Range("E16").Value = "  -4.40%  "

# Row 17
$ws.Range("E17").Value = "  -1.67%  "

# Row 18
Set-TextCell $ws.Range("D18") "18.25"
$ws.Range("E18").Value = "  -6.68%  "

# Row 19
$ws.Range("D19").Value = "66.668.83"
$ws.Range("E19").Value = "  -3.22%  "

# Row 20
Set-TextCell $ws.Range("D20") "12.08"
$ws.Range("E20").Value = "  -7.00%  "

# Row 21
$ws.Range("E21").Value = "  -8.41%  "

# Row 22
Set-TextCell $ws.Range("D22") "389.23"
$ws.Range("E22").Value = "  -5.85%  "

# Row 23
Set-TextCell $ws.Range("D23") "4.21"
$ws.Range("E23").Value = "  -8.32%  "

# Row 24
Set-TextCell $ws.Range("D24") "85.37"
$ws.Range("E24").Value = "  -5.28%  "

# Row 25
Set-TextCell $ws.Range("D25") "11.15"
$ws.Range("E25").Value = "  -1.48%  "

# Row 26
$ws.Range("E26").Value = "  -6.04%  "

# Row 27
Set-TextCell $ws.Range("D27") "12.40"
$ws.Range("E27").Value = "  -5.35%  "

# Row 28
Set-TextCell $ws.Range("D28") "6.05"
$ws.Range("E28").Value = "  +0.30%  "

# Row 29
Set-TextCell $ws.Range("D29") "3.57"
$ws.Range("E29").Value = "  -6.69%  "

# Row 30
Set-TextCell $ws.Range("D30") "8.84"
$ws.Range("E30").Value = "  -9.60%  "

# Row 31
Set-TextCell $ws.Range("D31") "7.53"
$ws.Range("E31").Value = "  -3.95%  "

# Row 32
Set-TextCell $ws.Range("D32") "30.87"
$ws.Range("E32").Value = "  -6.34%  "

# Row 33
Set-TextCell $ws.Range("D33") "628.59"
$ws.Range("E33").Value = "  -2.18%  "

# Row 34
Set-TextCell $ws.Range("D34") "12.14"
$ws.Range("E34").Value = "  -5.16%  "

# Row 35
$ws.Range("E35").Value = "  -7.79%  "

# Row 36
Set-TextCell $ws.Range("D36") "63.28"
$ws.Range("E36").Value = "  -6.45%  "

# Row 37
Set-TextCell $ws.Range("D37") "41.71"
$ws.Range("E37").Value = "  -9.85%  "

# Row 38
Set-TextCell $ws.Range("D38") "0.400"
$ws.Range("E38").Value = "  -4.25%  "

# Row 39
Set-TextCell $ws.Range("D39") "1.01"
$ws.Range("E39").Value = "  +0.50%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0752"
$ws.Range("E40").Value = "  -10.12%  "

# Row 41
$ws.Range("E41").Value = "  -5.75%  "

# Row 42
$ws.Range("D42").Value = "3.108.38"
$ws.Range("E42").Value = "  +6.01%  "

# Row 43
Set-TextCell $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
Set-TextCell $ws.Range("D44") "2.92"
$ws.Range("E44").Value = "  -5.16%  "

# Row 45
Set-TextCell $ws.Range("D45") "2.62"
$ws.Range("E45").Value = "  -0.62%  "

# Row 46
$ws.Range("E46").Value = "  -8.75%  "

# Row 47
Set-TextCell $ws.Range("D47") "3.14"
$ws.Range("E47").Value = "  +0.65%  "

# Row 48
$ws.Range("E48").Value = "  -7.36%  "

# Row 49
Set-TextCell $ws.Range("D49") "139.17"
$ws.Range("E49").Value = "  -4.72%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell $ws.Range("D50") "8.44"
$ws.Range("E50").Value = "  -9.27%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D51") "2.75"
$ws.Range("E51").Value = "  -1.30%  "
